$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) First (empty) paragraph in the document body: add <w:rFonts w:hint="cs"/>
#    as the first child of the run-properties on the paragraph mark, ahead of
#    the existing <w:highlight/> and <w:rtl/>.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$xmlParaHint = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p w:rsidR="00C82D58" w:rsidRDefault="00C82D58" w:rsidP="00EE36EE">' +
              '<w:pPr>' +
                '<w:rPr>' +
                  '<w:rFonts w:hint="cs"/>' +
                  '<w:highlight w:val="yellow"/>' +
                  '<w:rtl/>' +
                '</w:rPr>' +
              '</w:pPr>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'
$p1.Range.InsertXML($xmlParaHint)

# ---------------------------------------------------------------------------
# 2) Flip the "done?" flag of the first data row ("מוזיקת רקע, וכפתורים")
#    from "לא" to "כן" -- whole-word, first occurrence only.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("לא", $true, $true, $false, $false, $false, $true, 1, `
                         $false, "כן", 1) | Out-Null

# ---------------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark: it currently sits after the "כן" run in the
#    "הראה סיסמה" row; it must instead sit after the "כן" run we just wrote
#    in the "מוזיקת רקע, וכפתורים" row. The Bookmarks collection's own
#    Add/Delete calls don't round-trip reliably here, so the move is done by
#    rewriting each paragraph's OOXML in place (content otherwise unchanged).
# ---------------------------------------------------------------------------
$tbl = $d.Tables(1)

$cellWithNewBookmark = $tbl.Cell(2, 6)
$xmlAddBookmark = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p w:rsidR="005F585E" w:rsidRDefault="005F585E" w:rsidP="00160432">' +
              '<w:pPr>' +
                '<w:jc w:val="center"/>' +
                '<w:rPr>' +
                  '<w:rtl/>' +
                '</w:rPr>' +
              '</w:pPr>' +
              '<w:r>' +
                '<w:rPr>' +
                  '<w:rFonts w:hint="cs"/>' +
                  '<w:rtl/>' +
                '</w:rPr>' +
                '<w:t>כן</w:t>' +
              '</w:r>' +
              '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
              '<w:bookmarkEnd w:id="0"/>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'
$cellWithNewBookmark.Range.InsertXML($xmlAddBookmark)

$tbl = $d.Tables(1)
$cellWithOldBookmark = $tbl.Cell(6, 6)
$xmlRemoveBookmark = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p w:rsidR="005F585E" w:rsidRPr="00B70007" w:rsidRDefault="008117FC" w:rsidP="00160432">' +
              '<w:pPr>' +
                '<w:jc w:val="center"/>' +
                '<w:rPr>' +
                  '<w:color w:val="FFFFFF" w:themeColor="background1"/>' +
                  '<w:rtl/>' +
                '</w:rPr>' +
              '</w:pPr>' +
              '<w:r>' +
                '<w:rPr>' +
                  '<w:rFonts w:hint="cs"/>' +
                  '<w:color w:val="FFFFFF" w:themeColor="background1"/>' +
                  '<w:rtl/>' +
                '</w:rPr>' +
                '<w:t>כן</w:t>' +
              '</w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'
$cellWithOldBookmark.Range.InsertXML($xmlRemoveBookmark)

Write-Output "done"
